$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the day's refreshed
# coinranking.com snapshot (price + 1h volume-change columns), plus the
# Hedera / InternetComputer(DFINITY) rank swap (rows 33 & 34).
$updates = @{
    'D2' = '41.940.96'
    'E2' = '  +1.22%  '
    'D3' = '2.223.53'
    'E3' = '  +1.03%  '
    'D5' = '251.01'
    'E5' = '  -1.21%  '
    'E6' = '  -0.54%  '
    'D7' = '68.18'
    'E7' = '  -1.67%  '
    'E8' = '  +0.10%  '
    'D9' = '0.637'
    'E9' = '  +7.57%  '
    'D10' = '39.63'
    'E10' = '  +3.09%  '
    'D11' = '59.85'
    'E11' = '  +2.67%  '
    'D12' = '0.0941'
    'E12' = '  -1.21%  '
    'D13' = '7.10'
    'E13' = '  -1.99%  '
    'E14' = '  -0.08%  '
    'D15' = '2.555.25'
    'E15' = '  +1.18%  '
    'D16' = '14.70'
    'E16' = '  -0.85%  '
    'D17' = '0.876'
    'E17' = '  -1.12%  '
    'D18' = '2.217.83'
    'E18' = '  +1.81%  '
    'D19' = '41.864.25'
    'E19' = '  +1.31%  '
    'D20' = '0.0₃0964'
    'E20' = '  +0.80%  '
    'E21' = '  -0.26%  '
    'D22' = '72.82'
    'E22' = '  +1.04%  '
    'D23' = '232.68'
    'E23' = '  -0.21%  '
    'D24' = '2.08'
    'E24' = '  -0.46%  '
    'D25' = '3.91'
    'E25' = '  +0.82%  '
    'E26' = '  -5.01%  '
    'E27' = '  +0.13%  '
    'E28' = '  -4.64%  '
    'E29' = '  -1.59%  '
    'E30' = '  -1.80%  '
    'D31' = '167.13'
    'E31' = '  -1.95%  '
    'D32' = '20.46'
    'E32' = '  -1.15%  '
    'B33' = 'Hedera'
    'C33' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D33' = '0.0800'
    'E33' = '  +9.01%  '
    'B34' = 'InternetComputer(DFINITY)'
    'C34' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D34' = '5.96'
    'E34' = '  +6.63%  '
    'E35' = '  -0.75%  '
    'E36' = '  -0.52%  '
    'D37' = '4.65'
    'E37' = '  -0.38%  '
    'D38' = '4.12'
    'E38' = '  +2.14%  '
    'D39' = '25.39'
    'E39' = '  -3.19%  '
    'D40' = '0.0307'
    'E40' = '  +2.03%  '
    'E41' = '  +0.29%  '
    'D42' = '12.20'
    'E42' = '  +0.77%  '
    'E43' = '  -2.44%  '
    'D44' = '5.09'
    'E44' = '  +1.93%  '
    'D45' = '62.32'
    'E45' = '  -3.22%  '
    'D46' = '0.202'
    'E46' = '  -1.64%  '
    'D47' = '8.61'
    'E47' = '  -1.30%  '
    'D48' = '0.100'
    'E48' = '  -0.80%  '
    'D49' = '1.00'
    'E49' = '  -0.31%  '
    'E50' = '  +0.74%  '
    'D51' = '4.36'
    'E51' = '  +0.61%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Columns D/E hold free-text price/percentage strings (e.g. "41.940.96",
    # "0.100", "  +1.22%  "); force text format first so Excel's COM layer
    # doesn't silently reinterpret them as numbers and drop formatting like
    # trailing zeros.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
